$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Apio" was reported. It slots in at row 98
# (the table is not date-sorted), pushing the existing rows 98-120 down to
# 99-121.
$ws.Rows.Item(98).Insert()

# Populate the newly-inserted row 98 with the new record's data.
$ws.Cells.Item(98, 1).Value = 5
$ws.Cells.Item(98, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(98, 3).Value = "Maule"
$ws.Cells.Item(98, 4).Value = 44476
$ws.Cells.Item(98, 5).Value = 7
$ws.Cells.Item(98, 6).Value = 100112017
$ws.Cells.Item(98, 7).Value = "Apio"
$ws.Cells.Item(98, 8).Value = "Americana (o)"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 700
$ws.Cells.Item(98, 11).Value = 7000
$ws.Cells.Item(98, 12).Value = 7000
$ws.Cells.Item(98, 13).Value = 7000
$ws.Cells.Item(98, 14).Value = "`$/docena de matas"
$ws.Cells.Item(98, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(98, 16).Value = 1167
$ws.Cells.Item(98, 17).Value = 6
$ws.Cells.Item(98, 18).Value = "Hortaliza"
